$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert a new column before column B -----------------------------
# This shifts existing columns B,C,D (filename/name/effect) to C,D,E and
# shifts the cell contents + number formats along with it. Hyperlinks are
# NOT relocated automatically by this engine, so they are rebuilt below.
$ws.Columns("B").Insert()

# --- 2. Populate the new "default" column --------------------------------
# Write the "yes" values before the header "default" value so the shared
# string table gets them in the same order as in the target workbook.
$ws.Range("B2").Value = "yes"
$ws.Range("B3").Value = "yes"
$ws.Range("B4").Value = "yes"
$ws.Range("B5").Value = "yes"
$ws.Range("B1").Value = "default"

# --- 3. Rebuild the hyperlinks on the (now shifted) column E -------------
# This engine does not relocate hyperlink anchors when a column is
# inserted, so the previous D2:D5 hyperlinks are dropped and re-created on
# E2:E5. Hyperlinks.Add() stamps its TextToDisplay into both the cell AND
# the OOXML "display" attribute, but only the cell value can be changed
# afterwards (the "display" attribute stays), so: add with the URL as the
# display text, then restore each cell's real caption text, then restore
# the hyperlink cell style.
$disp1 = "https://www.pexels.com/es-es/foto/flor-flora-floracion-fondo-de-pantalla-de-flores-2187618/"
$disp2 = "https://www.pexels.com/es-es/foto/dalia-flor-flora-floracion-863963/"
$disp3 = "https://www.pexels.com/es-es/foto/flor-flora-floracion-fondo-de-pantalla-gratis-2039606/"
$disp4 = "https://www.pexels.com/es-es/foto/estambre-flor-flora-floracion-2375010/"

$text5 = $ws.Range("E5").Text
$text2 = $ws.Range("E2").Text
$text3 = $ws.Range("E3").Text
$text4 = $ws.Range("E4").Text

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E5"), $disp1, "", "", $disp1)
$ws.Range("E5").Value = $text5
$ws.Range("E5").Style = "Hipervínculo"

$ws.Hyperlinks.Add($ws.Range("E2"), $disp2, "", "", $disp2)
$ws.Range("E2").Value = $text2
$ws.Range("E2").Style = "Hipervínculo"

$ws.Hyperlinks.Add($ws.Range("E3"), $disp3, "", "", $disp3)
$ws.Range("E3").Value = $text3
$ws.Range("E3").Style = "Hipervínculo"

$ws.Hyperlinks.Add($ws.Range("E4"), $disp4, "", "", $disp4)
$ws.Range("E4").Value = $text4
$ws.Range("E4").Style = "Hipervínculo"

# --- 4. Column widths ------------------------------------------------------
# Column A keeps its existing width. The old column B/C/D widths now live on
# C/D (34.42578125 / 18.28515625) which Insert() already preserved. Give the
# new column B a narrow, auto-fit-like width for the short yes/no values.
$ws.Columns("B").ColumnWidth = 6.6

# --- 5. Selection ----------------------------------------------------------
$ws.Range("B12").Select()
